$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fill in the previously incomplete "Supplementary Table S3.1" row (old row 5) ---
# Write D before B to match the original authoring order (and resulting shared-string order).
$ws.Range("D5").Value2 = "Prepared. To be added in github."
$ws.Range("B5").Value2 = "Online Supplementary Material"

# --- Insert two new rows at the very top. ---
# Row 1 becomes the new title row, row 2 becomes a blank spacer row,
# and all previously existing rows (old 1-6) shift down to new rows 3-8.
$ws.Rows.Item(1).Insert()
$ws.Rows.Item(1).Insert()

# --- Title row ---
$ws.Range("A1").Value2 = "All Figures, Tables, Files (Main and Supplementary) of Chapter 3 in order of mention in manuscript."
$ws.Range("A1").Font.Bold = $false
$ws.Range("A1").Font.Size = 14
$ws.Rows.Item(1).RowHeight = 18.5

# --- Row 8 previously held the lone "Table 3.3" cell; replace it entirely ---
$ws.Range("A8").Value2 = "Supplementary Table S3.2"
$ws.Range("C8").Value2 = "Number of duplication and loss events for each gene family and comparison Cteno-first vs Sponge-first scenarios"
$ws.Range("B8").Value2 = "Online Supplementary Material"
$ws.Range("D8").Value2 = "Prepared. To be added in github."

# --- New row 9, only column D filled in ---
$ws.Range("D9").Value2 = "Prepared. To be added in github."

# --- Column width adjustments ---
$ws.Columns.Item(2).ColumnWidth = 28.08984375
$ws.Columns.Item(4).ColumnWidth = 29.1796875

# --- Selection matches final state ---
$ws.Range("D7:D9").Select()
